$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.082.00'
$ws.Range('E2').Value = '  +1.41%  '
$ws.Range('D3').Value = '1.892.11'
$ws.Range('E3').Value = '  +1.38%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.014'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +1.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '337.21'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.013'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.94%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4747'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.35%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3968'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.79%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '47.28'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08057'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.48%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.024'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.17%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.08'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.42%  '
$ws.Range('D13').Value = '1.900.53'
$ws.Range('E13').Value = '  +0.38%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.049'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.00%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.244'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.74%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.015'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.20%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '88.76'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.50%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06785'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.44%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.00001055'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.83%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.12'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.00%  '
$ws.Range('E21').Value = '  +0.94%  '
$ws.Range('D22').Value = '28.058.95'
$ws.Range('E22').Value = '  +1.27%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.549'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.22%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.06'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.87%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.350'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.68%  '
$ws.Range('D26').Value = '2.117.16'
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '160.99'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.43%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.12'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.122'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.70%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.565'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.33%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '122.18'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.37%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9846'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.81%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09593'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.28%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.644'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.32%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.376'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.30%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.378'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.64%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06104'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.65%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02260'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.35%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.208'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.20%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.250'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.52%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.012'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.87%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6003'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.42%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1901'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.80%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.40'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.72%  '
$ws.Range('E45').Value = '  +1.25%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5686'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.15%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '12.28'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.98%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.941'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.40%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.375'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.31%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06852'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.14%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '112.66'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.05%  '
